# Chat_Feature_Log.xlsx - "added events functionality stable"
#
# This script brings the Feature Log up to date:
#  - Clarifies the "View Resource - Link" backlog item with a question
#  - Un-hides / tidies a couple of entries
#  - Marks the "Events" feature as now InProgress (and highlights it green,
#    the sheet's existing convention for the row currently being worked on)
#  - Removes the highlight from "Filter on resource page" (no longer the
#    active focus) by resetting it back to a plain white fill
#  - Reworks the old "Group Chat - Display message..." idea into a clearer
#    "Improve System updates vs users chat" entry
#  - Adds three newly captured feature ideas to the bottom of the list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green highlight colour already used elsewhere on the sheet (FF92D050)
$Green = 5296274

# --- Resources section -----------------------------------------------

# B23: clarify whether this resource type is still required
$ws.Range("B23").Value = "View Resource - Link - is this really required?"

# --- Features section ---------------------------------------------------

# Row 31 ("Filter on resource page") loses its green highlight -> plain
# white fill (the sheet's "un-highlighted" background colour)
$ws.Rows(31).ClearFormats()
$ws.Range("A31").Interior.ThemeColor = 2
$ws.Range("B31").Interior.ThemeColor = 2
$ws.Range("D31").Interior.ThemeColor = 2

# Row 33 ("Events") is now in progress and gets the green highlight
$ws.Range("D33").Value = "InProgress"
$ws.Range("A33").Interior.Color = $Green
$ws.Range("B33").Interior.Color = $Green
$ws.Range("D33").Interior.Color = $Green

# B34: sharpen the wording of this feature idea
$ws.Range("B34").Value = "Group Chat - Improve System updates vs users chat (keep them separate)"

# New feature ideas captured at the bottom of the list
$ws.Range("A35").Value = "Features"
$ws.Range("B35").Value = "Joined Group List screen separate to current group screen"
$ws.Range("D35").Value = "Pending"

$ws.Range("A36").Value = "Features"
$ws.Range("B36").Value = "Registration ask for Country"
$ws.Range("D36").Value = "Pending"

$ws.Range("A37").Value = "Features"
$ws.Range("B37").Value = "Groups to be country specific"
$ws.Range("D37").Value = "Pending"
